$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New order rows (SKU, Name, Quantity, Cost Per, Total Cost) to append
# starting at row 3.
$data = @(
    @("124440",  "8 Grain",                             "4", "68.52",  "274.08"),
    @("108900",  "Seeds - Poppy",                        "1", "103.12", "103.12"),
    @("146200",  "Sugar - Donut",                        "2", "97.76",  "195.52"),
    @("4100860", "Boston Coffee Cake - Cinnamon Walnut", "3", "57.91",  "173.73"),
    @("4100870", "Boston Coffee Cake - Blueberry",       "4", "57.91",  "231.64"),
    @("4100760", "Boston Coffee Cake - Apple Cinnamon",  "3", "55.80",  "167.40")
)

$firstRow = 3
$lastRow = $firstRow + $data.Count - 1
$fullRange = $ws.Range("A$firstRow" + ":E$lastRow")

# Force text storage so numeric-looking values (SKU, qty, costs) are
# written as strings instead of being coerced into numbers.
$fullRange.NumberFormat = "@"

$r = $firstRow
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}

# Drop the temporary text format so the new cells stay unstyled, in
# keeping with the rest of the sheet.
$fullRange.ClearFormats()
